# Daily attendance processing - 2025-10-15 06:57:26
# Normalize the "Recorded By" (column G) cell values so that the
# "System"/"system" recorder is listed first in the comma separated list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "system, System, backup@backdoor.com") {
        $cell.Value2 = "System, system, backup@backdoor.com"
    }
}
